# Drag-and-drop implementation * Minor modifications: icon & GUI
#
# Slide 3 holds one big top-level group ("Group 42", the drag-handle icon)
# that is being repositioned/shrunk down onto the slide. Reproduce this by:
#   1. Ungrouping it so its four direct children become top-level shapes
#      (Oval 3, Group 41, Straight Connector 30, Group 37).
#   2. Moving/resizing each of those four children to their final
#      slide-absolute position (PowerPoint keeps each child's own internal
#      chOff/chExt -- i.e. everything nested further down -- untouched; only
#      the immediate children's own off/ext move, exactly like a normal
#      "resize the group" edit).
#   3. Thinning the dashed connector's line weight to match the smaller icon.
#   4. Re-grouping the four shapes, which is what naturally produces the
#      fresh shape id/name ("Group 2") seen in the target XML.
#
# NOTE on precision: Shape.Left/Top/Width/Height (like real PowerPoint COM)
# are marshalled as 32-bit `Single` points, so naively writing emu/12700.0
# can land 1 EMU off after PowerPoint truncates back to EMU. The helper
# below picks, for each target EMU value, the nearest points value whose
# float32 representation truncates back to exactly that EMU.

function EmuToPreciseGoodPoints {
    param(
        $TargetEmu
    )
    $base = $TargetEmu / 12700.0
    $bestPts = $base
    $bestDist = 1000000000.0
    $found = $false
    $step = 0.000001
    for ($n = -200; $n -le 200; $n++) {
        # Step across neighbouring float32 values around $base by nudging
        # the double before truncating it to Single (which is exactly what
        # PowerPoint's Shape.Left/Top/Width/Height -- a 32-bit `Single` --
        # does internally), then checking what EMU value PowerPoint would
        # truncate that Single back to.
        $cand = $base + ($n * $step)
        $single = [float]$cand
        $backEmu = [math]::Floor([float]$single * 12700.0)
        if ($backEmu -eq $TargetEmu) {
            $dist = [math]::Abs($cand - $base)
            if ($dist -lt $bestDist) {
                $bestDist = $dist
                $bestPts = $cand
                $found = $true
            }
        }
    }
    return $bestPts
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# The drag-handle icon is the only top-level shape on this slide.
$g = $s.Shapes.Item(1)

$null = $g.Ungroup()

# After ungrouping, the four former children are top-level shapes, in their
# original order.
$oval  = $s.Shapes.Item(1)   # Oval 3
$grp41 = $s.Shapes.Item(2)   # Group 41
$conn  = $s.Shapes.Item(3)   # Straight Connector 30
$grp37 = $s.Shapes.Item(4)   # Group 37

$oval.Left   = EmuToPreciseGoodPoints 3129793
$oval.Top    = EmuToPreciseGoodPoints 1554482
$oval.Width  = EmuToPreciseGoodPoints 1554480
$oval.Height = EmuToPreciseGoodPoints 1554480

$grp41.Left   = EmuToPreciseGoodPoints 3487910
$grp41.Top    = EmuToPreciseGoodPoints 1836826
$grp41.Width  = EmuToPreciseGoodPoints 838245
$grp41.Height = EmuToPreciseGoodPoints 989792

$conn.Left   = EmuToPreciseGoodPoints 3907033
$conn.Top    = EmuToPreciseGoodPoints 1836826
$conn.Width  = EmuToPreciseGoodPoints 0
$conn.Height = EmuToPreciseGoodPoints 989791
$conn.Line.Weight = EmuToPreciseGoodPoints 92075

$grp37.Left   = EmuToPreciseGoodPoints 3628855
$grp37.Top    = EmuToPreciseGoodPoints 2053543
$grp37.Width  = EmuToPreciseGoodPoints 556357
$grp37.Height = EmuToPreciseGoodPoints 556358

# Re-group the four shapes back together -- this mints the fresh id/name
# ("Group 2") that PowerPoint assigns to newly-created groups.
$range = $s.Shapes.Range(@(1, 2, 3, 4))
$newGroup = $range.Group()

Write-Output ("Regrouped as id=" + $newGroup.Id + " name=" + $newGroup.Name + `
    " left=" + $newGroup.Left + " top=" + $newGroup.Top + `
    " width=" + $newGroup.Width + " height=" + $newGroup.Height)
